# Apply updates to "VisualSA and Datalink Buffers" inputs
$wb = $excel.ActiveWorkbook

# --- "runs" sheet updates (tab literally named "runs" holds Start/Stop Run ID + Sim Time) ---
$runsTab = $wb.Worksheets.Item("runs")
$runsTab.Range("B2").Value = 1
$runsTab.Range("B3").Value = 250
$runsTab.Range("B2").Select()

# --- "params" sheet updates (tab literally named "params" holds the run configuration table) ---
$paramsTab = $wb.Worksheets.Item("params")
$paramsTab.Range("D2").Value = 0
$paramsTab.Range("E2").Value = 5
$paramsTab.Range("E3").Select()

$paramsTab.Activate()
